$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "B-30", 150, 20, 15.35860333660845, 0),
    @(2, "B-39", 150, 20, 13.1984467046723, 0),
    @(3, "B-48", 150, 20, 11.67311575022249, 0),
    @(4, "B-63", 150, 20, 12.21065474546192, 0),
    @(5, "B-13", 150, 20, 11.66818684759706, 0),
    @(6, "B-56", 150, 20, 11.72986436997653, 0),
    @(7, "B-179", 150, 20, 11.5072466797277, 0),
    @(8, "B-180", 150, 20, 11.50356061438601, 0),
    @(9, "B-21", 300, 20, 36.50135474911401, 0)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # Copy the formatting (style) of A2 -> column A of the new row,
    # matching the bold/border/center style used for the index column.
    $ws.Range("A2").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $r++
}

$excel.CutCopyMode = 0
